# ---------------------------------------------------------------------------
# Applies the two changes captured in the commit:
#   1. The table on slide 16 gets a different built-in table style applied
#      (Table.ApplyStyle - mirrors clicking a different style in the Table
#      Design gallery).
#   2. The deck's theme colour scheme (the one actually driving the slide
#      master / layouts / slides) is switched from the custom "Integral"
#      palette back to the stock "Office" palette (Design tab -> Colors ->
#      "Office"). PowerPoint stores colours as the legacy VBA RGB() encoding
#      (0xBBGGRR) on ThemeColorScheme.Colors(i).RGB, so the hex triplets from
#      the target theme are byte-swapped below.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 16 -------------------------------------------
$s = $p.Slides.Item(16)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{A63E1542-73D7-45FA-B64F-FE264586F382}")
    }
}

# --- 2. Theme colour scheme: Integral -> Office ----------------------------
# Index : Office colour : RGB() value (0xBBGGRR)
#   1 dk1       000000 -> 0
#   2 lt1       FFFFFF -> 16777215
#   3 dk2       44546A -> 6968388
#   4 lt2       E7E6E6 -> 15132391
#   5 accent1   5B9BD5 -> 13998939
#   6 accent2   ED7D31 -> 3243501
#   7 accent3   A5A5A5 -> 10855845
#   8 accent4   FFC000 -> 49407
#   9 accent5   4472C4 -> 12874308
#  10 accent6   70AD47 -> 4697456
#  11 hlink     0563C1 -> 12673797
#  12 folHlink  954F72 -> 7491477
$officeRgb = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

$master = $p.SlideMaster
$themeColors = $master.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = $officeRgb[$i - 1]
}
